# Proyecto Final - Historias de Usuario: Modificacion 2.0
# Renumber the HU (Historia de Usuario) items starting at HU#3 (the
# "REQUEST" story) so that the newly inserted HU#3 ("modelo" story,
# already present at row 26) keeps its number and every following
# story/task shifts up by one: HU#3->HU#4, HU#4->HU#5, HU#5->HU#6,
# HU#6->HU#7 (title) while its DATASOURCE tasks become HU#8_T1/T2,
# HU#7->HU#9 and HU#8->HU#10. Column A text only; B/C/D stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "HU#4 Como desarrollador quiero que el proyecto tenga validación de los parámetros de entrada."
$ws.Range("A30").Value = "HU#4_T1 Creación de los tests correspondientes a la etapa REQUEST."
$ws.Range("A31").Value = "HU#4_T2 Creación de las clases correspondientes a la etapa REQUEST con su correspondiente implementación."

$ws.Range("A32").Value = "HU#5 Como Product Owner quiero que el proyecto disponga de controladores para controlar las comunicaciones entre clases de la etapa REQUEST y las de la etapa SERVICE."
$ws.Range("A33").Value = "HU#5_T1 Creación de los tests correspondientes a la etapa CONTROLLER."
$ws.Range("A34").Value = "HU#5_T2 Creación de las clases correspondientes a la etapa CONTROLLER con su correspondiente implementación."

$ws.Range("A35").Value = "HU#6 Como desarrollador quiero que el proyecto tenga unas clases que provean el servicio e implementen la lógica de negocio."
$ws.Range("A36").Value = "HU#6_T1 Creación de los tests correspondientes a la etapa SERVICE."
$ws.Range("A37").Value = "HU#6_T2 Creación de las clases correspondientes a la etapa SERVICE con su correspondiente implementación."

$ws.Range("A38").Value = "HU#7 Como manager del proyecto quiero que el proyecto disponga de unas interfaces que se encarguen de dar respuesta a las peticiones de la etapa SERVICE."

$ws.Range("A39").Value = "HU#8_T1 Creación de los tests correspondientes a la etapa DATASOURCE."
$ws.Range("A40").Value = "HU#8_T2 Creación de las clases correspondientes a la etapa DATASOURCE con su correspondiente implementación."

$ws.Range("A41").Value = "HU#9 Como Product Owner quiero que el proyecto tenga unas clases para que se encarguen de realizar llamadas a la API y recibir su respuesta con su correspondiente tratamiento."
$ws.Range("A42").Value = "HU#9_T1 Creación de los tests correspondientes a las clases que llamen a la API."
$ws.Range("A43").Value = "HU#9_T2 Creación de las clases que realicen llamadas a la API de coinlore."

$ws.Range("A44").Value = "HU#10 Como desarrollador quiero que el proyecto almacene los datos obtenidos de la API en caché."
$ws.Range("A45").Value = "HU#10_T1 Creación de los tests correspondientes a las clases que se encarguen del almacenamiento de los datos en caché."
$ws.Range("A46").Value = "HU#10_T2 Creación de las clases que se encarguen del almacenamiento de los datos en caché."

# Move the selection/active view to match the author's final cursor
# position (was A26, now A45) and scroll the sheet up so row 4 is the
# top-visible row (was row 17).
$ws.Range("A45").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
